$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.691.82"
$ws.Range("E2").Value = "  +2.35%  "

# Row 3
$ws.Range("D3").Value = "2.609.08"
$ws.Range("E3").Value = "  +1.14%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.87%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "

# Row 8
$ws.Range("E8").Value = "  +0.62%  "

# Row 9
$ws.Range("D9").Value = "2.635.12"
$ws.Range("E9").Value = "  +1.85%  "

# Row 10
$ws.Range("E10").Value = "  -2.20%  "

# Row 11
$ws.Range("E11").Value = "  +1.81%  "

# Row 12
$ws.Range("E12").Value = "  -5.91%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.365"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.91%  "

# Row 14
$ws.Range("D14").Value = "3.076.86"
$ws.Range("E14").Value = "  +1.40%  "

# Row 15
$ws.Range("D15").Value = "60.672.94"
$ws.Range("E15").Value = "  +2.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000142"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.70%  "

# Row 18
$ws.Range("D18").Value = "2.466.57"
$ws.Range("E18").Value = "  -4.64%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.43%  "

# Row 20
$ws.Range("E20").Value = "  +1.98%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.49%  "

# Row 23
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("E24").Value = "  +12.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.28%  "

# Row 27
$ws.Range("E27").Value = "  +0.16%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.46%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0795"
$ws.Range("E29").Value = "  +2.63%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.03%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.27%  "

# Row 32
$ws.Range("E32").Value = "  -0.14%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.55%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.80%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.981"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.99%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.38%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.56%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.844"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.73%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "297.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.84%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.88%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0986"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.607"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.54%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0547"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0242"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.16%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.32%  "
